$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WEEK1")
$ws.Activate()

# Insert two new blank rows before row 5 (shifts existing rows 5-13 down to 7-15)
$ws.Rows("5:6").Insert()

# New row 5: team meeting entry, part 1 (no description cell)
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Group"
$ws.Range("D5").Value = "29/7/2020"

# New row 6: team meeting entry, part 2 (with description)
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Group"
$ws.Range("C6").Value = "Fixed some problems related to database and namespace. Sharing some ideas "
$ws.Range("D6").Value = "29/7/2020"

# Match the row heights used for similar "Group" meeting rows
$ws.Rows(5).RowHeight = 36.6
$ws.Rows(6).RowHeight = 36.6

# Update the active selection
$ws.Range("C6").Select()
